$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'" + "37.817.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'" + "  -0.83%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'" + "2.028.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'" + "  -1.60%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'" + "  -0.12%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'" + "227.24"
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'" + "0.613"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'" + "  -0.46%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = "'" + "  +1.30%  "
$ws.Range("E7").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'" + "0.383"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'" + "  -1.00%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'" + "0.0814"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'" + "  +0.52%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("E11").Value = "'" + "  -0.20%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'" + "14.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'" + "  -0.63%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'" + "2.329.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'" + "  -1.65%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'" + "20.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'" + "  +1.14%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("E15").Value = "'" + "  +0.40%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'" + "5.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'" + "  -1.83%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'" + "2.054.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'" + "  -0.55%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'" + "37.724.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'" + "  -0.78%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'" + "6.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'" + "  -1.76%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'" + "69.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'" + "  -0.24%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'" + "0.0₃0823"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'" + "  -1.05%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'" + "224.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'" + "  -0.08%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'" + "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'" + "  -0.03%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("E24").Value = "'" + "  -2.60%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("E25").Value = "'" + "  -2.02%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("B26").Value = "'" + "Cosmos"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'" + "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'" + "9.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'" + "  -0.86%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("B27").Value = "'" + "Monero"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'" + "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'" + "165.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'" + "  -0.73%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'" + "0.127"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'" + "  -3.88%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'" + "18.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'" + "  -0.75%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = "'" + "  -4.81%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = "'" + "  +1.04%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("E32").Value = "'" + "  -2.92%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "'" + "  +5.20%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("E34").Value = "'" + "  -3.03%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("E35").Value = "'" + "  -2.18%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'" + "6.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'" + "  +5.42%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = "'" + "  -4.35%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("E38").Value = "'" + "  -2.92%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("E39").Value = "'" + "  +0.00%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'" + "1.528.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'" + "  +3.22%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("E41").Value = "'" + "  -0.95%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("B42").Value = "'" + "FTXToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'" + "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'" + "4.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'" + "  +7.71%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("B43").Value = "'" + "Aave"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'" + "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'" + "96.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'" + "  -1.96%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("B44").Value = "'" + "InjectiveProtocol"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'" + "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'" + "16.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'" + "  -0.44%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("E45").Value = "'" + "  -0.53%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'" + "0.0915"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'" + "  -2.99%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("E47").Value = "'" + "  -1.85%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("E48").Value = "'" + "  -1.70%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'" + "  -0.40%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = "'" + "  +0.19%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'" + "2.218.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'" + "  -1.63%  "
$ws.Range("E51").Style = "Normal"
